# "Version 1." -> "Version 2.", re-shaped into the exact run layout the
# author's (re)commit produced:
#   "Versi" | "on" | " 2" | <bookmarkStart/End _GoBack> | "."
#
# Word normally fuses adjacent runs that share identical formatting, so a
# plain Find/Replace collapses everything back into a single <w:r>. To force
# a genuine run boundary (with no visible formatting change) we briefly drop
# a same-formatted bookmark at the split point and remove it again; the
# split survives because the two pieces were separated by the bookmark when
# they were created.

$d = $word.ActiveDocument

# --- 1. Split "Version" -> "Versi" + "on" at offset 5 -----------------
$splitPoint = $d.Range(5, 5)
$d.Bookmarks.Add("TmpRunSplit", $splitPoint) | Out-Null
$d.Bookmarks("TmpRunSplit").Delete()

# --- 2. "1" -> "2" (offset 8-9, the digit right after the space) ------
$d.Range(8, 9).Text = "2"

# --- 3. Move the trailing "." so it lands after the _GoBack bookmark --
# Drop it from the " 2." run ...
$d.Range(9, 10).Text = ""
# ... then re-insert it past the (now-trailing) _GoBack bookmark so it
# becomes its own run following the bookmark, matching the target XML.
$d.Range(9, 9).InsertAfter(".")
